# "Generate Report for Archive"
#
# The localization-status report was regenerated: the status that used to
# read "Ready for handoff" is now "In Translation" everywhere it appears
# (the Overview roll-up sheet, per-locale columns E/F, and the "Status"
# column on each per-locale detail sheet). Because the new status text is
# shorter than the old one, the report generator also re-sized the
# "Status"/locale columns to fit the new text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn (col E) / de-de (col F) status cells ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

# Shrink columns E:F to fit the new (shorter) status text.
$overview.Range("E1:F1").EntireColumn.ColumnWidth = 12.5

# --- Per-locale detail sheets: "Status" column (col C) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C1").EntireColumn.ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C1").EntireColumn.ColumnWidth = 12.5
